$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 19.02.2022 01:15"

# Update row 5 (Makro) values: shift current price to old price column, set new price,
# compute delta as text, and store the check timestamp as text in E5
$ws.Range("B5").Value = 36.7
$ws.Range("C5").Value = 36.5

# Force D5 to stay as literal text "+0.2" (otherwise Excel would parse it as a number),
# then drop back to the default (unstyled) cell format.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "+0.2"
$ws.Range("D5").Style = "Normal"

# E5 switches from a formatted date-serial number to plain text, losing its date style.
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2022-02-19 01:15:14"
$ws.Range("E5").Style = "Normal"
